# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Column G ("K") values are recalculated/rewritten with the new strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 4
    3  = 6
    4  = 4
    5  = 4
    6  = 8
    7  = 3
    8  = 8
    9  = 2
    10 = 7
    11 = 5
    12 = 2
    13 = 7
    14 = 4
    15 = 6
    16 = 3
    17 = 1
    18 = 2
    19 = 4
    20 = 2
    21 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
